$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 970
$wsExhibit.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202403/fmbmIP421710756195423.jpeg"
$wsExhibit.Range("F3").Value = 1941
$wsExhibit.Range("F4").Value = 428

# Sheet "全部类型" (all types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 970
$wsAll.Range("I4").Value = "//i1.hdslb.com/bfs/openplatform/202403/fmbmIP421710756195423.jpeg"
$wsAll.Range("F5").Value = 1941
$wsAll.Range("F6").Value = 428
